$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row at 197, pushing existing rows 197:211 down to 198:212
$ws.Rows.Item(197).Insert()

# Populate the newly inserted row 197 with the new record
$ws.Cells.Item(197, 1).Value = 5
$ws.Cells.Item(197, 2).Value = "Macroferia Regional de Talca"
$ws.Cells.Item(197, 3).Value = "Maule"
$ws.Cells.Item(197, 4).Value = 45267
$ws.Cells.Item(197, 5).Value = 7
$ws.Cells.Item(197, 6).Value = "Fruta"
$ws.Cells.Item(197, 7).Value = 100108
$ws.Cells.Item(197, 8).Value = "Tropicales y subtropicales"
$ws.Cells.Item(197, 9).Value = 100108002
$ws.Cells.Item(197, 10).Value = "Mango"
$ws.Cells.Item(197, 11).Value = "Sin especificar"
$ws.Cells.Item(197, 12).Value = "Primera"
$ws.Cells.Item(197, 13).Value = 248
$ws.Cells.Item(197, 14).Value = 11000
$ws.Cells.Item(197, 15).Value = 11000
$ws.Cells.Item(197, 16).Value = 11000
$ws.Cells.Item(197, 17).Value = "`$/bandeja 4 kilos"
$ws.Cells.Item(197, 18).Value = "Perú"
$ws.Cells.Item(197, 19).Value = 2750
$ws.Cells.Item(197, 20).Value = 4
